$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("letter", "手紙|てがみ"),
    @("singer", "歌手|かしゅ"),
    @("hand", "手|て"),
    @("sign language", "手話|しゅわ"),
    @("good at", "上手な|じょうずな"),
    @("paper", "紙|かみ"),
    @("Japanese paper", "和紙|わし"),
    @("cover page", "表紙|ひょうし"),
    @("origami", "折り紙|おりがみ"),
    @("to like", "好きな|すきな"),
    @("to love", "大好きな|だいすきな"),
    @("goodwill", "好意|こうい"),
    @("liking; taste", "好み|このみ"),
    @("favorite food", "好物|こうぶつ"),
    @("near; nearby", "近く|ちかく"),
    @("neighborhood", "近所|きんじょ"),
    @("recently", "最近|さいきん"),
    @("the Middle and Near East", "中近東|ちゅうきんとう"),
    @("cheerful; bright", "明るい|あかるい"),
    @("tomorrow", "明日|あした"),
    @("explanation", "説明|せつめい"),
    @("invention", "発明|はつめい"),
    @("civilization", "文明|ぶんめい"),
    @("hospital", "病院|びょういん"),
    @("illness", "病気|びょうき"),
    @("serious illness", "重病|じゅうびょう"),
    @("sudden illness", "急病|きゅうびょう"),
    @("graduate school", "大学院|だいがくいん"),
    @("beauty parlor", "美容院|びよういん"),
    @("movie", "映画|えいが"),
    @("movie theater", "映画館|えいがかん"),
    @("to be reflected", "映る|うつる"),
    @("painter", "画家|がか"),
    @("plan", "計画|けいかく"),
    @("comic", "漫画|まんが"),
    @("to sing", "歌う|うたう"),
    @("song", "歌|うた"),
    @("singer", "歌手|かしゅ"),
    @("national anthem", "国歌|こっか"),
    @("Kabuki", "歌舞伎|かぶき"),
    @("lyrics", "歌詞|かし"),
    @("Kawaguchi City", "川口市|かわぐちし"),
    @("city hall", "市役所|しやくしょ"),
    @("mayor", "市長|しちょう"),
    @("market", "市場|いちば"),
    @("various places", "いろいろな所|いろいろなところ"),
    @("neighborhood", "近所|きんじょ"),
    @("kitchen", "台所|だいどころ"),
    @("address", "住所|じゅうしょ"),
    @("to study", "勉強する|べんきょうする"),
    @("to try hard", "勉める|つとめる"),
    @("diligent", "勤勉な|きんべんな"),
    @("to study", "勉強する|べんきょうする"),
    @("strong", "強い|つよい"),
    @("obstinate", "強情な|ごうじょうな"),
    @("robbery", "強盗|ごうとう"),
    @("powerful", "強力な|きょうりょくな"),
    @("famous", "有名な|ゆうめいな"),
    @("toll; fee", "有料|ゆうりょう"),
    @("to exist", "有る|ある"),
    @("talented", "有能な|ゆうのうな"),
    @("travel", "旅行|りょこう"),
    @("Japanese inn", "旅館|りょかん"),
    @("traveling alone", "一人旅|ひとりたび"),
    @("passport", "旅券|りょけん"),
)

$startRow = 100
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    [void]($ws.Cells.Item($row, 1).Value = $data[$i][0])
    [void]($ws.Cells.Item($row, 2).Value = $data[$i][1])
}
